# bug fixes on test 5 generator
# Updates the Write Latency "min" (O) and "max" (P) columns for each
# benchmark row (3-38) on Sheet1 with corrected values from the fixed
# test-5 result generator. Values are plain numeric / "NNNk" style text
# (matching the rest of the table), so they are written through a
# scratch helper cell with a leading apostrophe and pasted in as
# values-only; this keeps the destination cells text-typed without
# disturbing their existing cell style (s="1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("O3", "587"),
    @("P3", "701124"),
    @("O4", "985"),
    @("P4", "2283.8k"),
    @("O5", "578"),
    @("P5", "833566"),
    @("O6", "569"),
    @("P6", "453980"),
    @("O7", "479"),
    @("P7", "20263k"),
    @("O8", "587"),
    @("P8", "1193.0k"),
    @("O9", "500"),
    @("P9", "783755"),
    @("O10", "392"),
    @("P10", "826038"),
    @("O11", "643"),
    @("P11", "891430"),
    @("O12", "635"),
    @("P12", "812513"),
    @("O13", "897"),
    @("P13", "893861"),
    @("O14", "905"),
    @("P14", "2941.4k"),
    @("O15", "551"),
    @("P15", "646099"),
    @("O16", "627"),
    @("P16", "810721"),
    @("O17", "402"),
    @("P17", "837399"),
    @("O18", "646"),
    @("P18", "461274"),
    @("O19", "623"),
    @("P19", "1705.6k"),
    @("O20", "767"),
    @("P20", "1187.3k"),
    @("O21", "699"),
    @("P21", "829431"),
    @("O22", "652"),
    @("P22", "669866"),
    @("O23", "455"),
    @("P23", "881095"),
    @("O24", "702"),
    @("P24", "736755"),
    @("O25", "477"),
    @("P25", "492913"),
    @("O26", "545"),
    @("P26", "783086"),
    @("O27", "452"),
    @("P27", "433936"),
    @("O28", "497"),
    @("P28", "759788"),
    @("O29", "608"),
    @("P29", "1453.4k"),
    @("O30", "694"),
    @("P30", "1649.9k"),
    @("O31", "496"),
    @("P31", "839150"),
    @("O32", "539"),
    @("P32", "894701"),
    @("O33", "501"),
    @("P33", "526823"),
    @("O34", "488"),
    @("P34", "783065"),
    @("O35", "515"),
    @("P35", "521805"),
    @("O36", "647"),
    @("P36", "817326"),
    @("O37", "533"),
    @("P37", "1312.8k"),
    @("O38", "828"),
    @("P38", "481918")
)

# Scratch cell well outside the A1:X38 table used to stage each text value.
$helper = $ws.Range("Z1")

foreach ($update in $updates) {
    $cellRef = $update[0]
    $newValue = $update[1]
    $helper.Value = "'" + $newValue
    $helper.Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}

$helper.Clear() | Out-Null
$excel.CutCopyMode = $false
